$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''60.772.18'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = '''2.988.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -6.06%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.05%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''561.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -5.24%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''126.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -7.11%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.08%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''2.984.82'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -6.09%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  -3.44%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '''  -6.62%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''5.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -2.94%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = '''  -4.76%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '''  -7.22%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''32.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -6.79%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = '''  +0.33%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''3.487.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -5.99%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''60.869.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -3.57%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''2.991.64'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -5.95%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''6.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -7.52%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''430.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -6.70%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''13.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -6.83%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''0.657'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -7.73%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -7.66%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''12.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -4.11%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''78.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -6.01%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.09%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.32%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''2.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -8.00%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''7.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -9.32%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -8.35%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''25.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -8.25%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''5.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -12.45%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''0.0928'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -10.39%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''2.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -5.86%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.952'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -8.82%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''5.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -5.20%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''49.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -2.86%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.0₃0663'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -7.05%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.0356'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -8.77%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''7.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -4.90%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = '''Bittensor'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''372.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -8.96%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = '''Kaspa'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''0.106'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -5.80%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''2.671.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -4.67%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''2.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -8.97%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.05%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''TheGraph'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''0.233'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -8.05%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''Monero'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''119.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -3.24%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''33.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -3.64%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''1.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -9.39%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  -5.40%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''23.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -9.46%  '
$ws.Range('E51').Style = 'Normal'
